$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.262.16"
$ws.Range("E2").Value = "  +0.14%  "

# Row 3
$ws.Range("D3").Value = "2.542.63"
$ws.Range("E3").Value = "  -2.71%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "'588.65"
$ws.Range("E5").Value = "  -0.60%  "

# Row 6
$ws.Range("D6").Value = "'172.58"
$ws.Range("E6").Value = "  +3.89%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").Value = "'0.529"
$ws.Range("E8").Value = "  -0.50%  "

# Row 9
$ws.Range("D9").Value = "2.540.81"
$ws.Range("E9").Value = "  -2.76%  "

# Row 10
$ws.Range("D10").Value = "'0.143"
$ws.Range("E10").Value = "  +3.25%  "

# Row 11
$ws.Range("D11").Value = "'0.161"
$ws.Range("E11").Value = "  +0.20%  "

# Row 12
$ws.Range("D12").Value = "'0.353"
$ws.Range("E12").Value = "  -2.62%  "

# Row 13
$ws.Range("D13").Value = "'5.16"
$ws.Range("E13").Value = "  -0.80%  "

# Row 14
$ws.Range("D14").Value = "'27.10"
$ws.Range("E14").Value = "  -1.29%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000179"
$ws.Range("E15").Value = "  -0.76%  "

# Row 16
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.968.80"
$ws.Range("E16").Value = "  -3.92%  "

# Row 17
$ws.Range("D17").Value = "67.086.37"
$ws.Range("E17").Value = "  -0.13%  "

# Row 18
$ws.Range("D18").Value = "2.553.69"
$ws.Range("E18").Value = "  -1.94%  "

# Row 19
$ws.Range("D19").Value = "'7.99"
$ws.Range("E19").Value = "  +1.25%  "

# Row 20
$ws.Range("D20").Value = "'11.50"
$ws.Range("E20").Value = "  -2.95%  "

# Row 21
$ws.Range("D21").Value = "'352.52"
$ws.Range("E21").Value = "  -0.81%  "

# Row 22
$ws.Range("D22").Value = "'4.24"
$ws.Range("E22").Value = "  -1.41%  "

# Row 23
$ws.Range("D23").Value = "'4.73"
$ws.Range("E23").Value = "  +1.94%  "

# Row 24
$ws.Range("D24").Value = "'1.99"
$ws.Range("E24").Value = "  +3.45%  "

# Row 25
$ws.Range("E25").Value = "  +0.02%  "

# Row 26
$ws.Range("D26").Value = "'70.21"
$ws.Range("E26").Value = "  +1.29%  "

# Row 27
$ws.Range("D27").Value = "'10.00"
$ws.Range("E27").Value = "  -2.62%  "

# Row 28
$ws.Range("E28").Value = "  -0.01%  "

# Row 29
$ws.Range("D29").Value = "2.644.40"
$ws.Range("E29").Value = "  -4.05%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0998"
$ws.Range("E30").Value = "  -0.47%  "

# Row 31
$ws.Range("D31").Value = "'534.94"
$ws.Range("E31").Value = "  -1.09%  "

# Row 32
$ws.Range("D32").Value = "'8.33"
$ws.Range("E32").Value = "  +5.21%  "

# Row 33
$ws.Range("D33").Value = "'1.34"
$ws.Range("E33").Value = "  -0.55%  "

# Row 34
$ws.Range("D34").Value = "'1.86"
$ws.Range("E34").Value = "  -1.67%  "

# Row 35
$ws.Range("D35").Value = "'0.132"
$ws.Range("E35").Value = "  -2.92%  "

# Row 36
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  -0.12%  "

# Row 37
$ws.Range("D37").Value = "'1.49"
$ws.Range("E37").Value = "  -0.31%  "

# Row 38
$ws.Range("D38").Value = "'157.19"
$ws.Range("E38").Value = "  +0.30%  "

# Row 39
$ws.Range("D39").Value = "'18.79"
$ws.Range("E39").Value = "  -0.90%  "

# Row 40
$ws.Range("D40").Value = "'18.43"
$ws.Range("E40").Value = "  +1.09%  "

# Row 41
$ws.Range("D41").Value = "'0.358"
$ws.Range("E41").Value = "  -1.83%  "

# Row 42
$ws.Range("D42").Value = "'1.82"
$ws.Range("E42").Value = "  +1.27%  "

# Row 43
$ws.Range("D43").Value = "'5.15"
$ws.Range("E43").Value = "  -0.59%  "

# Row 44
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.54"
$ws.Range("E44").Value = "  +5.20%  "

# Row 45
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.02%  "

# Row 46
$ws.Range("D46").Value = "'151.85"
$ws.Range("E46").Value = "  -0.13%  "

# Row 47
$ws.Range("D47").Value = "'0.565"
$ws.Range("E47").Value = "  -2.10%  "

# Row 48
$ws.Range("D48").Value = "'3.73"
$ws.Range("E48").Value = "  -1.31%  "

# Row 49
$ws.Range("D49").Value = "'1.75"
$ws.Range("E49").Value = "  +3.20%  "

# Row 50
$ws.Range("D50").Value = "0.0₆0277"
$ws.Range("E50").Value = "  -8.11%  "

# Row 51
$ws.Range("D51").Value = "'0.0761"
$ws.Range("E51").Value = "  -0.98%  "
